$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 42, shifting existing data
# (rows 42-55) down to rows 44-57.
$ws.Rows.Item(42).Insert()
$ws.Rows.Item(42).Insert()

# New row 42: "Primera" quality record for the 2022-07-13 market date.
$ws.Cells.Item(42, 1).Value = 11
$ws.Cells.Item(42, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(42, 3).Value = "Bíobío"
$ws.Cells.Item(42, 4).Value = 44755
$ws.Cells.Item(42, 5).Value = 8
$ws.Cells.Item(42, 6).Value = 100112043
$ws.Cells.Item(42, 7).Value = "Pepino dulce"
$ws.Cells.Item(42, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 100
$ws.Cells.Item(42, 11).Value = 14000
$ws.Cells.Item(42, 12).Value = 15000
$ws.Cells.Item(42, 13).Value = 14500
$ws.Cells.Item(42, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(42, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(42, 16).Value = 806
$ws.Cells.Item(42, 17).Value = 18
$ws.Cells.Item(42, 18).Value = "Hortaliza"

# New row 43: "Segunda" quality record for the 2022-07-13 market date.
$ws.Cells.Item(43, 1).Value = 11
$ws.Cells.Item(43, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(43, 3).Value = "Bíobío"
$ws.Cells.Item(43, 4).Value = 44755
$ws.Cells.Item(43, 5).Value = 8
$ws.Cells.Item(43, 6).Value = 100112043
$ws.Cells.Item(43, 7).Value = "Pepino dulce"
$ws.Cells.Item(43, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(43, 9).Value = "Segunda"
$ws.Cells.Item(43, 10).Value = 50
$ws.Cells.Item(43, 11).Value = 11000
$ws.Cells.Item(43, 12).Value = 11000
$ws.Cells.Item(43, 13).Value = 11000
$ws.Cells.Item(43, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(43, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(43, 16).Value = 611
$ws.Cells.Item(43, 17).Value = 18
$ws.Cells.Item(43, 18).Value = "Hortaliza"
